$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 795.6286
$ws.Range("I15").Value = 795.6286
$ws.Range("K15").Value = 2386.8858
$ws.Range("M15").Value = -2217.8858

$ws.Range("H21").Value = 12007.5
$ws.Range("J21").Value = 11998
$ws.Range("L21").Value = 11998
$ws.Range("N21").Value = -12934

$ws.Range("H23").Value = 12007.5
$ws.Range("J23").Value = 11998
$ws.Range("L23").Value = 11998
$ws.Range("N23").Value = -12466

$ws.Range("H29").Value = 83333460
$ws.Range("I29").Value = 83333460
$ws.Range("K29").Value = 250000380
$ws.Range("M29").Value = -250000099

$ws.Range("H32").Value = 3058.2
$ws.Range("I32").Value = 2895
$ws.Range("J32").Value = 3167
$ws.Range("K32").Value = 2895
$ws.Range("L32").Value = 3167
$ws.Range("M32").Value = -2569
$ws.Range("N32").Value = -3819

$ws.Range("H33").Value = 763.9677
$ws.Range("I33").Value = 139.95454
$ws.Range("K33").Value = 139.95454
$ws.Range("M33").Value = 89.04545999999999

$ws.Range("H43").Value = 9349.75
$ws.Range("I43").Value = 1900
$ws.Range("K43").Value = 1900
$ws.Range("M43").Value = -1831

$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H100").Value = 4596
$ws.Range("J100").Value = 5450.8887
$ws.Range("L100").Value = 5450.8887
$ws.Range("N100").Value = -6532.8887

$ws.Range("H129").Value = 1182.2354
$ws.Range("I129").Value = 679.7778
$ws.Range("K129").Value = 2039.3334
$ws.Range("M129").Value = 2960.6666

$ws.Range("H132").Value = 2411.739
$ws.Range("I132").Value = 2051.3157
$ws.Range("J132").Value = 4123.75
$ws.Range("K132").Value = 6153.9471
$ws.Range("L132").Value = 12371.25
$ws.Range("M132").Value = -3623.9471
$ws.Range("N132").Value = -17431.25

$ws.Range("H137").Value = 11229.167
$ws.Range("I137").Value = 13512.23
$ws.Range("J137").Value = 5293.2
$ws.Range("K137").Value = 40536.69
$ws.Range("L137").Value = 15879.6
$ws.Range("M137").Value = -37986.69
$ws.Range("N137").Value = -20979.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1765.4333
$ws.Range("I2").Value = 772.82355
$ws.Range("J2").Value = 3063.4614
$ws.Range("K2").Value = 772.82355
$ws.Range("L2").Value = 3063.4614
$ws.Range("M2").Value = -659.82355
$ws.Range("N2").Value = -3289.4614

$ws.Range("H37").Value = 9627.200000000001
$ws.Range("I37").Value = 2034
$ws.Range("K37").Value = 2034
$ws.Range("M37").Value = -1761

$ws.Range("H74").Value = 2138.423
$ws.Range("I74").Value = 2177.56
$ws.Range("J74").Value = 1160
$ws.Range("K74").Value = 2177.56
$ws.Range("L74").Value = 1160
$ws.Range("M74").Value = -1303.56
$ws.Range("N74").Value = -2908

$ws.Range("H77").Value = 2138.423
$ws.Range("I77").Value = 2177.56
$ws.Range("J77").Value = 1160
$ws.Range("K77").Value = 10887.8
$ws.Range("L77").Value = 5800
$ws.Range("M77").Value = -6519.799999999999
$ws.Range("N77").Value = -14536

$ws.Range("H116").Value = 1765.4333
$ws.Range("I116").Value = 772.82355
$ws.Range("J116").Value = 3063.4614
$ws.Range("K116").Value = 772.82355
$ws.Range("L116").Value = 3063.4614
$ws.Range("M116").Value = 1521.17645
$ws.Range("N116").Value = -7651.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1765.4333
$ws.Range("I3").Value = 772.82355
$ws.Range("J3").Value = 3063.4614
$ws.Range("K3").Value = 772.82355
$ws.Range("L3").Value = 3063.4614
$ws.Range("M3").Value = -658.82355
$ws.Range("N3").Value = -3291.4614

$ws.Range("H94").Value = 2208.4348
$ws.Range("I94").Value = 2055.2778
$ws.Range("K94").Value = 2055.2778
$ws.Range("M94").Value = -1604.2778

$ws.Range("H99").Value = 3084.32
$ws.Range("I99").Value = 2094.3125
$ws.Range("J99").Value = 4844.3335
$ws.Range("K99").Value = 2094.3125
$ws.Range("L99").Value = 4844.3335
$ws.Range("M99").Value = -596.3125
$ws.Range("N99").Value = -7840.3335

$ws.Range("H107").Value = 10000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 10000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4334.6294
$ws.Range("I31").Value = 1892
$ws.Range("J31").Value = 15082.2
$ws.Range("K31").Value = 1892
$ws.Range("L31").Value = 15082.2
$ws.Range("M31").Value = -1597
$ws.Range("N31").Value = -15672.2

$ws.Range("H34").Value = 4334.6294
$ws.Range("I34").Value = 1892
$ws.Range("J34").Value = 15082.2
$ws.Range("K34").Value = 1892
$ws.Range("L34").Value = 15082.2
$ws.Range("M34").Value = -1690
$ws.Range("N34").Value = -15486.2

$ws.Range("H56").Value = 49747
$ws.Range("I56").Value = 24495
$ws.Range("J56").Value = 74999
$ws.Range("K56").Value = 24495
$ws.Range("L56").Value = 74999
$ws.Range("M56").Value = -23650
$ws.Range("N56").Value = -76689

$ws.Range("H107").Value = 1136.0264
$ws.Range("I107").Value = 351.42856
$ws.Range("K107").Value = 351.42856
$ws.Range("M107").Value = 1568.57144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 133.55556
$ws.Range("I33").Value = 121
$ws.Range("J33").Value = 143.6
$ws.Range("K33").Value = 726
$ws.Range("L33").Value = 861.5999999999999
$ws.Range("M33").Value = -443
$ws.Range("N33").Value = -1427.6

$ws.Range("H98").Value = 1128.8334
$ws.Range("J98").Value = 1383
$ws.Range("L98").Value = 4149
$ws.Range("N98").Value = -7145

$ws.Range("H109").Value = 219097.8
$ws.Range("I109").Value = 258872.25
$ws.Range("J109").Value = 60000
$ws.Range("K109").Value = 776616.75
$ws.Range("L109").Value = 180000
$ws.Range("M109").Value = -775576.75
$ws.Range("N109").Value = -182080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7914.1816
$ws.Range("I80").Value = 10117.75
$ws.Range("J80").Value = 5269.9
$ws.Range("K80").Value = 10117.75
$ws.Range("L80").Value = 5269.9
$ws.Range("M80").Value = -9119.75
$ws.Range("N80").Value = -7265.9

$ws.Range("H83").Value = 7914.1816
$ws.Range("I83").Value = 10117.75
$ws.Range("J83").Value = 5269.9
$ws.Range("K83").Value = 50588.75
$ws.Range("L83").Value = 26349.5
$ws.Range("M83").Value = -45596.75
$ws.Range("N83").Value = -36333.5

$ws.Range("H97").Value = 1257.125
$ws.Range("J97").Value = 1134.2858
$ws.Range("L97").Value = 1134.2858
$ws.Range("N97").Value = -2126.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 4650
$ws.Range("I17").Value = 2944.4443
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 2944.4443
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = -2774.4443
$ws.Range("N17").Value = -20340

$ws.Range("H46").Value = 2499.875
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H93").Value = 4390.5713
$ws.Range("I93").Value = 3423.647
$ws.Range("K93").Value = 3423.647
$ws.Range("M93").Value = -2175.647
